# Update "想去人数" (want-to-go count) values in column F for the
# 展览 (sheet1) and 全部类型 (sheet4) worksheets, reflecting the refreshed
# data snapshot ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value for 展览 sheet
$exhibitUpdates = @{
    3  = 1366
    7  = 1211
    8  = 1524
    9  = 155
    11 = 1079
    12 = 438
    13 = 104
    16 = 85
    18 = 6041
    20 = 5875
    21 = 9839
    24 = 179
    25 = 270
    26 = 492
    27 = 161
    28 = 143
    29 = 4376
    30 = 366
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value for 全部类型 sheet (same events, different row numbers)
$allUpdates = @{
    5  = 1366
    10 = 1211
    12 = 1524
    14 = 155
    15 = 1079
    17 = 438
    18 = 104
    22 = 85
    24 = 6041
    26 = 5875
    27 = 9839
    31 = 179
    32 = 270
    34 = 492
    38 = 161
    39 = 143
    40 = 4376
    46 = 366
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
